$d = $word.ActiveDocument

# Remove the leading empty paragraph (widowControl/pBdr/spacing/ind/jc formatting,
# no text) that precedes the "LAMPIRAN A" heading paragraph.
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Delete()

# Remove the now-orphan bookmark ("_heading=h.gjdgxs") wrapping the start of the
# "LAMPIRAN A" heading paragraph.
$bm = $d.Bookmarks("_heading=h.gjdgxs")
$bm.Delete()
